# Generate Report for Archive
# - Update status text "Ready for handoff" -> "In Translation" everywhere it appears
# - Narrow the date/status related columns to match new content width

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    foreach ($cell in $used.Cells) {
        if ("Ready for handoff" -eq $cell.Value2) {
            $cell.Value2 = "In Translation"
        }
    }
}

# Target column width (OOXML "width" attribute) is 13.4101845877511.
# The ColumnWidth COM property is specified in "characters" and this runtime
# (matching genuine Excel pixel-snapping behaviour) stores width internally
# in whole-pixel steps, i.e. raw = ColumnWidth + 5/6. Back solving gives the
# characters value that reproduces the desired stored width as closely as
# the pixel grid allows.
$newColumnWidth = 13.4101845877511 - (5 / 6)

# Overview sheet: columns E (zh-cn) and F (de-de) status columns
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Columns.Item(5).ColumnWidth = $newColumnWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newColumnWidth

# zh-cn sheet: column C (Status)
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Columns.Item(3).ColumnWidth = $newColumnWidth

# de-de sheet: column C (Status)
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Columns.Item(3).ColumnWidth = $newColumnWidth
